# Historial de abonos y cambios Sistema Joyeria Violeta
# Commit: "Se arregla reimpresión de pagos, y se limpian productos cuando
#          cambian el tipo de venta a taller"
#
# Net effect on Sheet1:
#  - Row 14/15 were both blank placeholder rows under the "Modificaciones"
#    list; row 14 now gets a new entry (barcode printing feature) and the
#    redundant blank row 15 is removed (everything below shifts up by one).
#  - Three blank filler rows near the bottom of the "Abonos" list are
#    trimmed down to just one blank row (two rows removed), so everything
#    from the "Total de abonos"/"Pendiente" rows down shifts up by two
#    more. The SUM()/subtraction formulas re-anchor automatically when the
#    rows are deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the extra blank row above "total acomulado" ------------------
$ws.Rows.Item(15).Delete()

# --- trim the 3 blank filler rows before the totals down to just 1 -------
# (after the row-15 delete above, these blank rows now live at 29/30/31)
$ws.Rows.Item(29).Delete()
$ws.Rows.Item(29).Delete()

# --- fill in the new "Modificaciones" line item on (what is now) row 14 --
$ws.Range("B14").Value = "Impresión de código de barra de producto, Configuración de código de barras en impresora de oficina, Modificaciones de impresión de ticket de notas, Corte de Caja, Egresos, Ingresos y pagos"
$ws.Range("C14").Value = 2600
$ws.Rows.Item(14).RowHeight = 42.75
